$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:G values for rows 2-38 (regenerated sval data filtering save games)
$data = @(
    @(0.6545652718822623,0.3048912486333797,3.223369029078222,0.5333859586016987,0,4.716211508195562),
    @(0.0006075818656279264,0.002658071450198252,3.223369029078222,0.5333859586016987,1,3.760020640995746),
    @(1.445647641019636,1.626987699542094,0.7210945179870265,0.5333859586016987,1,4.327115817150455),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,0,6.15379541431027),
    @(1.445647641019636,1.626987699542094,18.71679738969934,0.5333859586016987,0,22.32281868886277),
    @(3.272327238179451,1.626987699542094,0.1496068669990043,0.5333859586016987,1,5.582307763322248),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(1.445647641019636,1.626987699542094,0.7210945179870265,0.5333859586016987,1,4.327115817150455),
    @(3.272327238179451,1.626987699542094,18.71679738969934,0.5333859586016987,1,24.14949828602258),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,18.71679738969934,0.5333859586016987,1,24.14949828602258),
    @(3.272327238179451,1.626987699542094,0.1496068669990043,0.5333859586016987,1,5.582307763322248),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,1,6.15379541431027),
    @(3.272327238179451,1.626987699542094,0.1496068669990043,0.5333859586016987,1,5.582307763322248),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,18.71679738969934,0.5333859586016987,1,24.14949828602258),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,1,6.15379541431027),
    @(0.6545652718822623,1.626987699542094,3.223369029078222,0.5333859586016987,1,6.038307959104277),
    @(0.6545652718822623,1.626987699542094,0.7210945179870265,0.5333859586016987,1,3.536033448013082),
    @(1.445647641019636,1.626987699542094,189.6080260415259,0.5333859586016987,1,193.2140473406893),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,18.71679738969934,0.5333859586016987,1,24.14949828602258),
    @(3.272327238179451,1.626987699542094,0.1496068669990043,0.5333859586016987,1,5.582307763322248),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,0,8.656069925401464),
    @(0.6545652718822623,1.626987699542094,3.223369029078222,0.5333859586016987,1,6.038307959104277),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,1,6.15379541431027),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,0,6.15379541431027),
    @(1.445647641019636,1.626987699542094,3.223369029078222,0.5333859586016987,1,6.82939032824165),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,0,6.15379541431027),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,13.86384647080068,1,19.48425592650926),
    @(3.272327238179451,1.626987699542094,18.71679738969934,13.86384647080068,1,37.47995879822157),
    @(3.272327238179451,1.626987699542094,0.7210945179870265,0.5333859586016987,0,6.15379541431027),
    @(3.272327238179451,1.626987699542094,3.223369029078222,0.5333859586016987,1,8.656069925401464),
    @(3.272327238179451,1.626987699542094,18.71679738969934,0.5333859586016987,1,24.14949828602258)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value2 = $vals[$j]
    }
}
